$wb = $excel.ActiveWorkbook

# Insert a new worksheet "MessagePage" right before "TeamMembers" (after "InvoicePage")
$before = $wb.Worksheets.Item("TeamMembers")
$ws = $wb.Worksheets.Add($before)
$ws.Name = "MessagePage"

$ws.Range("A1").Value = "Subject"
$ws.Range("B1").Value = "TestSubject"
$ws.Range("A2").Value = "Message"
$ws.Range("B2").Value = "TestMessage"

$ws.Activate()
